$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.422"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05635"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.424"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.362"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8122"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9170"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1442"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07496"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03089"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03112"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09344"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.555"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001586"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04770"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005791"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006389"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.004998"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.190"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3301"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1295"
$ws.Range("D26").Style = "Normal"
$ws.Range("D40").Value = "'0.04040"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006801"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002711"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007553"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005698"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
